$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.260.16"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.908.97"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +9.85%  "
$ws.Range("D6").Value = "'253.26"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'40.53"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "'0.363"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").Value = "'52.28"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'0.0764"
$ws.Range("E11").Value = "  +6.68%  "
$ws.Range("D12").Value = "'0.0988"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "2.185.70"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'12.70"
$ws.Range("E14").Value = "  +5.13%  "
$ws.Range("D15").Value = "'0.717"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "1.926.52"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'4.91"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "35.253.37"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'74.35"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "0.0₃0846"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "'243.27"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "'13.02"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "'5.07"
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.37"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").Value = "'166.72"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "'8.63"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "'18.71"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +4.63%  "
$ws.Range("D31").Value = "4.126.65"
$ws.Range("E31").Value = "  +19.40%  "
$ws.Range("D32").Value = "'4.34"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.00"
$ws.Range("E33").Value = "  +14.51%  "
$ws.Range("B34").Value = "TrustWalletToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D34").Value = "'1.63"
$ws.Range("E34").Value = "  +22.41%  "
$ws.Range("D35").Value = "'0.0582"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "'2.03"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0217"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.06"
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").Value = "'96.37"
$ws.Range("E42").Value = "  +6.84%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'0.0648"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").Value = "1.337.04"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "'2.77"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").Value = "'45.32"
$ws.Range("E50").Value = "  -5.84%  "
$ws.Range("D51").Value = "'12.04"
$ws.Range("E51").Value = "  +16.06%  "
